# Insert a new weekly price record for "Haba" (Macroferia Regional de Talca)
# at row 18, shifting the existing rows 18-61 down to 19-62. This mirrors a
# new week of data being added at the top of the chronological block while
# the rest of the historic rows (previously rows 18-61) move down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 18:61 down to 19:62, leaving row 18 blank (format/style of row 18
# -- notably the date style on column D -- is inherited from the row above).
$ws.Rows("18:18").Insert()

# Populate the newly freed row 18 with the new data point.
$ws.Cells.Item(18, 1).Value = 5
$ws.Cells.Item(18, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(18, 3).Value = "Maule"
$ws.Cells.Item(18, 4).Value = 44519
$ws.Cells.Item(18, 5).Value = 7
$ws.Cells.Item(18, 6).Value = 100112026
$ws.Cells.Item(18, 7).Value = "Haba"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 6000
$ws.Cells.Item(18, 12).Value = 6000
$ws.Cells.Item(18, 13).Value = 6000
$ws.Cells.Item(18, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(18, 15).Value = "Región del Maule"
$ws.Cells.Item(18, 16).Value = 240
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"
